$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Password column (E) for rows 7, 10, and 11 from "password" to "test"
$ws.Range("E7").Value = "test"
$ws.Range("E10").Value = "test"
$ws.Range("E11").Value = "test"
